# Generate Report for Handback
# Update the Correspond Handoff / Handback datetimes on the per-language
# sheets to reflect the newly generated report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-30 10:46:39"
$wsZhCn.Range("H2").Value = "2016-03-30 10:47:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-30 10:46:50"
$wsDeDe.Range("H2").Value = "2016-03-30 10:47:45"
